# Update "想去人数" (want-to-go count) figures that changed in the latest
# data refresh for the gh-pages generated output.
#
# These updates apply to both the "展览" sheet and the "全部类型" sheet,
# which mirror the same underlying exhibition data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 7710
    "F6"  = 40
    "F9"  = 5942
    "F12" = 29
    "F13" = 1811
    "F14" = 1326
    "F15" = 284
    "F16" = 584
    "F17" = 131
    "F18" = 5536
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
